$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row 33 with data, following the pattern of existing rows
$ws.Cells.Item(33, 1).Value = 10002
$ws.Cells.Item(33, 2).Value = 10032
$ws.Cells.Item(33, 3).Value = "eng"
$ws.Cells.Item(33, 4).Value = $true
$ws.Cells.Item(33, 5).Value = "superadmin"
$ws.Cells.Item(33, 6).Value = "now()"
$ws.Cells.Item(33, 7).Value = "now()"

# Update selection to match the target state
$ws.Range("B30").Select()
